$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.118.44'
$ws.Range("E2").Value = '  +2.79%  '
$ws.Range("D3").Value = '2.314.05'
$ws.Range("E3").Value = '  +2.57%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.32'
$ws.Range("E5").Value = '  +1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.36'
$ws.Range("E6").Value = '  +6.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.539'
$ws.Range("E7").Value = '  +2.66%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("E9").Value = '  +7.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.19'
$ws.Range("E10").Value = '  +4.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0818'
$ws.Range("E11").Value = '  +3.73%  '
$ws.Range("E12").Value = '  +0.88%  '
$ws.Range("E13").Value = '  +7.95%  '
$ws.Range("D14").Value = '2.670.41'
$ws.Range("E14").Value = '  +2.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.01'
$ws.Range("E15").Value = '  +4.84%  '
$ws.Range("D16").Value = '2.315.97'
$ws.Range("E16").Value = '  +2.51%  '
$ws.Range("E17").Value = '  +3.30%  '
$ws.Range("D18").Value = '43.052.23'
$ws.Range("E18").Value = '  +2.92%  '
$ws.Range("E19").Value = '  +1.96%  '
$ws.Range("E20").Value = '  +2.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.12'
$ws.Range("E21").Value = '  +3.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.56'
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.41'
$ws.Range("E23").Value = '  +1.94%  '
$ws.Range("E24").Value = '  +5.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.64'
$ws.Range("E25").Value = '  +3.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.71'
$ws.Range("E27").Value = '  +4.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.45'
$ws.Range("E28").Value = '  +3.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.67'
$ws.Range("E29").Value = '  +2.74%  '
$ws.Range("E30").Value = '  -0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '166.05'
$ws.Range("E31").Value = '  +3.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.35'
$ws.Range("E32").Value = '  +3.47%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.13'
$ws.Range("E34").Value = '  -1.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.94'
$ws.Range("E35").Value = '  +5.58%  '
$ws.Range("E36").Value = '  +1.40%  '
$ws.Range("E37").Value = '  +3.70%  '
$ws.Range("E38").Value = '  +1.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.85'
$ws.Range("E39").Value = '  +2.34%  '
$ws.Range("E40").Value = '  +2.22%  '
$ws.Range("E41").Value = '  +8.73%  '
$ws.Range("E42").Value = '  +1.18%  '
$ws.Range("E43").Value = '  +2.87%  '
$ws.Range("D44").Value = '1.976.28'
$ws.Range("E44").Value = '  +0.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.09'
$ws.Range("E45").Value = '  +2.20%  '
$ws.Range("E46").Value = '  +4.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.83'
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.95'
$ws.Range("E48").Value = '  +18.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.68'
$ws.Range("E49").Value = '  +5.64%  '
$ws.Range("D50").Value = '2.538.15'
$ws.Range("E50").Value = '  +2.30%  '
$ws.Range("E51").Value = '  +3.09%  '
